$wb = $excel.ActiveWorkbook

# Sheets "展览" (index 1) and "全部类型" (index 4) both contain the same
# two-event table and both need the new event inserted as the new row 2,
# pushing the existing "丽水·幻梦动漫嘉年华" entry down to row 3 (and
# renumbering its index from 1 to 2).
$sheetIndexes = @(1, 4)

foreach ($idx in $sheetIndexes) {
    $ws = $wb.Worksheets.Item($idx)

    # Insert a blank row above the current data row (row 2), shifting the
    # existing event down to row 3.
    $ws.Rows.Item(2).Insert()

    # Restore the thin border on the new A2 cell to match the other index
    # cells in column A (row insert only copies font/alignment, not border).
    $ws.Cells.Item(2, 1).Borders.LineStyle = 1

    # Fill in the new row 2 with the "丽水·CCAC动漫游戏嘉年华" event.
    $ws.Cells.Item(2, 1).Value = 1

    # Force the date-looking string to stay plain text (not get parsed into
    # a date serial number), then drop back to the default "Normal" style
    # so no extra number formatting lingers on the cell.
    $ws.Cells.Item(2, 2).Value = "'2024-11-24"
    $ws.Cells.Item(2, 2).Style = "Normal"

    $ws.Cells.Item(2, 3).Value = "丽水·CCAC动漫游戏嘉年华"
    $ws.Cells.Item(2, 4).Value = "中东路848号(解放街交汇) 飞达国际大酒店"
    $ws.Cells.Item(2, 5).Value = "2024.11.24 09:00-11.24 17:00"
    $ws.Cells.Item(2, 6).Value = 1
    $ws.Cells.Item(2, 7).Value = 29.9
    $ws.Cells.Item(2, 8).Value = "https://show.bilibili.com/platform/detail.html?id=93797"
    $ws.Cells.Item(2, 9).Value = "//i2.hdslb.com/bfs/openplatform/202410/zK2vYBx41729481548356.jpeg"

    # The event that used to be row 2 is now row 3; bump its index from 1 to 2.
    $ws.Cells.Item(3, 1).Value = 2
}
